# Update the Developer_Feedback sheet so the hyperlink display text points to
# "GitHub URL" labels instead of the old Dropbox-hosted .msg file paths, and
# update the footnote to reflect the new GitHub (public) access location.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Developer_Feedback")

# Jenkins rows (D5:D7)
$ws.Range("D5").Value = "Jenkins GitHub URL1"
$ws.Range("D6").Value = "Jenkins GitHub URL 2"
$ws.Range("D7").Value = "Jenkins GitHub URL 3"

# Wordpress rows (D8:D9)
$ws.Range("D8").Value = "Wordpress GitHub URL 1"
$ws.Range("D9").Value = "Wordpress GitHub URL 2"

# IdleonCompanion rows (D10:D11)
$ws.Range("D10").Value = "IdleonCompanion GitHub URL 1"
$ws.Range("D11").Value = "IdleonCompanion GitHub URL 2"

# OpenMW rows (D12:D13)
$ws.Range("D12").Value = "OpenMW GitHub URL 1"
$ws.Range("D13").Value = "OpenMW GitHub URL 2"

# Footnote
$ws.Range("A24").Value = "* - The above links for text and message files and are accessed from GitHub (public) location"
